$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rerun of the HW1 retrieval-eval code on a single sharded index lowered
# the scores for the "TF-IDF" (row 5) and "Okapi BM25" (row 6) models.
$ws.Range("B4").Value = 0.3
$ws.Range("B5").Value = 0.2962
$ws.Range("C5").Value = 0.42
$ws.Range("D5").Value = 0.3333
$ws.Range("B6").Value = 0.225
$ws.Range("C6").Value = 0.42
$ws.Range("D6").Value = 0.332

# The author's active cell/selection ended up on B6 when the file was saved.
$null = $ws.Range("B6").Select()
